$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price/volume updates (no coin name or link changes) ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.717.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.465.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.16%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.61"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.66"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.73%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.62"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.87%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.400"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.055.43"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.73"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.79%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.466.66"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.781.21"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.38"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.40"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.22"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.54"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.29%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.93"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.602.01"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000115"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.182"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.60"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.38%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.13"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.13"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.27%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.78"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.30"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.10"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.19"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +20.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.502.92"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0769"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.799"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.49"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.15"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.72"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.19"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.602.42"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.36"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.01%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.42%  "


# --- Row 37/38 swap: Monero <-> ImmutableX (with updated price/volume) ---
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.58"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.73%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "170.59"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.17%  "

# --- Row 50/51 swap: FirstDigitalUSD <-> dogwifhat (with updated price/volume) ---
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.16%  "

$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.09%  "
